$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DateOfGame is now recorded as free-form text instead of a date serial
$ws.Range("E2").Value = "2021-25-04"

# Update the game's host country (Brazil -> Turkey)
$ws.Range("C2").Value = "Turkey"

# Dynamic percentage of system accuracy
$ws.Range("K2").Value = 62

# Drop the one-off fill flag that K2 used to carry; falls back to the
# plain bordered style shared by the rest of the data row.
$ws.Range("K2").Interior.Pattern = -4142

# Column K is now a touch narrower (stored width 15 chars)
$ws.Columns.Item(11).ColumnWidth = 14.15

# Scroll/selection state as left by the author
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("J5").Select()
